$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue "D2" "38.130.96"
Set-TextValue "E2" "  +2.96%  "
Set-TextValue "D3" "2.060.73"
Set-TextValue "E3" "  +2.60%  "
Set-TextValue "E4" "  -0.32%  "
Set-TextValue "D5" "230.40"
Set-TextValue "E5" "  +1.98%  "
Set-TextValue "E6" "  +2.77%  "
Set-TextValue "D7" "58.31"
Set-TextValue "E7" "  +6.60%  "
Set-TextValue "E8" "  +0.00%  "
Set-TextValue "E9" "  +2.73%  "
Set-TextValue "D10" "0.0807"
Set-TextValue "E10" "  +2.62%  "
Set-TextValue "E11" "  -1.16%  "
Set-TextValue "D12" "2.364.89"
Set-TextValue "E12" "  +2.55%  "
Set-TextValue "D13" "14.63"
Set-TextValue "E13" "  +3.52%  "
Set-TextValue "E14" "  +2.30%  "
Set-TextValue "D15" "0.756"
Set-TextValue "E15" "  +2.38%  "
Set-TextValue "E16" "  +3.99%  "
Set-TextValue "D17" "2.058.11"
Set-TextValue "E17" "  +2.60%  "
Set-TextValue "D18" "38.022.11"
Set-TextValue "E18" "  +2.90%  "
Set-TextValue "D19" "6.18"
Set-TextValue "E19" "  +1.02%  "
Set-TextValue "D20" "69.74"
Set-TextValue "E20" "  +1.46%  "
Set-TextValue "E21" "  +1.71%  "
Set-TextValue "D22" "224.78"
Set-TextValue "E22" "  +0.75%  "
Set-TextValue "E23" "  -0.01%  "
Set-TextValue "E24" "  +1.22%  "
Set-TextValue "E25" "  +3.05%  "
Set-TextValue "D26" "9.34"
Set-TextValue "E26" "  +2.59%  "
Set-TextValue "D27" "165.95"
Set-TextValue "D28" "0.134"
Set-TextValue "E28" "  +8.15%  "
Set-TextValue "D29" "19.08"
Set-TextValue "E29" "  +2.14%  "
Set-TextValue "E30" "  +2.49%  "
Set-TextValue "E31" "  +1.88%  "
Set-TextValue "E32" "  +1.25%  "
Set-TextValue "D33" "4.62"
Set-TextValue "E33" "  +4.68%  "
Set-TextValue "E34" "  +1.10%  "
Set-TextValue "D35" "1.99"
Set-TextValue "E35" "  +7.48%  "
Set-TextValue "E36" "  +1.93%  "
Set-TextValue "E37" "  +13.57%  "
Set-TextValue "D38" "3.32"
Set-TextValue "E38" "  +5.67%  "
Set-TextValue "E39" "  -0.27%  "
Set-TextValue "D40" "98.49"
Set-TextValue "E40" "  +4.03%  "
Set-TextValue "E41" "  +1.13%  "
Set-TextValue "D42" "1.484.73"
Set-TextValue "E42" "  +0.71%  "
Set-TextValue "B43" "InjectiveProtocol"
Set-TextValue "C43" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D43" "16.88"
Set-TextValue "E43" "  +2.14%  "
Set-TextValue "B44" "Cronos"
Set-TextValue "C44" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D44" "0.0945"
Set-TextValue "E44" "  +2.81%  "
Set-TextValue "E45" "  +3.90%  "
Set-TextValue "E46" "  +0.24%  "
Set-TextValue "D47" "4.08"
Set-TextValue "E47" "  +16.76%  "
Set-TextValue "E48" "  +1.37%  "
Set-TextValue "D49" "2.97"
Set-TextValue "E49" "  +2.23%  "
Set-TextValue "E50" "  -0.92%  "
Set-TextValue "D51" "2.252.69"
Set-TextValue "E51" "  +2.70%  "
